$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database - rows 2 and 3 (Estonia, Furn/Home Furnishings)
foreach ($row in 2, 3) {
    $ws.Range("D$row").Value = -0.0989

    $ws.Range("G$row").Value = -0.002462686567164179
    $ws.Range("H$row").Value = -0.002462686567164179
    $ws.Range("I$row").Value = -0.01902985074626866
    $ws.Range("J$row").Value = -0.01902985074626866
    $ws.Range("K$row").Value = -0.272
    $ws.Range("L$row").Value = -0.02029850746268657

    $ws.Range("U$row").Value = 0.005
    $ws.Range("V$row").Value = 0.002016129032258064
    $ws.Range("W$row").Value = -0.1554285714285714
    $ws.Range("X$row").Value = 0.1584895902618158
    $ws.Range("Y$row").Value = -0.3139181616903872
    $ws.Range("Z$row").Value = 1.99345432906873
    $ws.Range("AA$row").Value = -0.03793513835168105
    $ws.Range("AB$row").Value = 0.06818970680525233
    $ws.Range("AC$row").Value = -0.1061248451569334
    $ws.Range("AD$row").Value = 5.26
    $ws.Range("AE$row").Value = 0
    $ws.Range("AF$row").Value = 5.26
    $ws.Range("AG$row").Value = 5.255
    $ws.Range("AH$row").Value = 0.6795865633074935
    $ws.Range("AI$row").Value = 0.7601156069364161
    $ws.Range("AJ$row").Value = 0.6793794440853265
    $ws.Range("AK$row").Value = 0.759942154736081
    $ws.Range("AL$row").Value = 0.244
    $ws.Range("AM$row").Value = 0.244
    $ws.Range("AN$row").Value = 10.84536082474227
    $ws.Range("AO$row").Value = -1.045081967213115
    $ws.Range("AP$row").Value = 10.83505154639175
    $ws.Range("AQ$row").Value = -1.045081967213115
}
